$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp shown at the top of the sheet
$ws.Range("A1").Value = "Datos actualizados a 9 de Abril de 2020 a las 01:22"

# --- Estados Unidos (row 4): refreshed case counts, no reordering ---
$ws.Range("B4").Value = 427079
$ws.Range("C4").Value = 26744
$ws.Range("D4").Value = 22314
$ws.Range("E4").Value = 390100
$ws.Range("G4").Value = 1824
$ws.Range("H4").Value = 14665

# --- Noruega / Irlanda swap (rows 24-25) ---
# Irlanda moves up to row 24 (its own totals carry over unchanged),
# Noruega drops to row 25 and receives refreshed totals.
$ws.Range("A24").Value = "Irlanda"
$ws.Range("B24").Value = 6074
$ws.Range("C24").Value = 365
$ws.Range("D24").Value = 25
$ws.Range("E24").Value = 5814
$ws.Range("F24").Value = 165
$ws.Range("G24").Value = 25
$ws.Range("H24").Value = 235

$ws.Range("A25").Value = "Noruega"
$ws.Range("B25").Value = 6042
$ws.Range("C25").Value = 0
$ws.Range("D25").Value = 32
$ws.Range("E25").Value = 5909
$ws.Range("F25").Value = 78
$ws.Range("G25").Value = 12
$ws.Range("H25").Value = 101

# --- Guayana Francesa moves ahead of Monaco / Liechtenstein / Aruba (rows 133-136) ---
# Guayana Francesa takes row 133 with refreshed totals; Monaco, Liechtenstein and
# Aruba each shift down one row, carrying their previous totals unchanged.
$ws.Range("A133").Value = "Guayana Francesa"
$ws.Range("B133").Value = 83
$ws.Range("C133").Value = 6
$ws.Range("D133").Value = 43
$ws.Range("E133").Value = 40
$ws.Range("F133").Value = 1
$ws.Range("G133").Value = 0
$ws.Range("H133").Value = 0

$ws.Range("A134").Value = "Monaco"
$ws.Range("B134").Value = 81
$ws.Range("C134").Value = 2
$ws.Range("D134").Value = 4
$ws.Range("E134").Value = 76
$ws.Range("F134").Value = 4
$ws.Range("G134").Value = 0
$ws.Range("H134").Value = 1

$ws.Range("A135").Value = "Liechtenstein"
$ws.Range("B135").Value = 78
$ws.Range("C135").Value = 0
$ws.Range("D135").Value = 55
$ws.Range("E135").Value = 22
$ws.Range("F135").Value = 0
$ws.Range("G135").Value = 0
$ws.Range("H135").Value = 1

$ws.Range("A136").Value = "Aruba"
$ws.Range("B136").Value = 77
$ws.Range("C136").Value = 3
$ws.Range("D136").Value = 14
$ws.Range("E136").Value = 63
$ws.Range("F136").Value = 0
$ws.Range("G136").Value = 0
$ws.Range("H136").Value = 0
